$d = $word.ActiveDocument

# Locate the heading paragraph "Implemented plausibility checks", then grab the
# (empty) paragraph right after it -- that is the placeholder that the new
# safety-requirement paragraphs replace.
$target = $null
$paras = $d.Paragraphs
$count = $paras.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $paras.Item($i)
    if ($para.Range.Text -match "Implemented plausibility checks") {
        $target = $paras.Item($i + 1)
        break
    }
}

$xml = @"
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>BrakePedalPressed</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>boolean</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>value</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> so plausibility checks are not needed.</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>ThrottlePedalPosition</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> needs to be in</w:t></w:r><w:r><w:t xml:space="preserve"> the</w:t></w:r><w:r><w:t xml:space="preserve"> range [0; 1] so if the position is not valid, a warning </w:t></w:r><w:r><w:t>flag is set</w:t></w:r><w:r><w:t>, the torque is set to 0 and vehicle switches to Neutral State.</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>AutomaticTransmissionSelectorState</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is seen as an integer value between 0 and 4. </w:t></w:r><w:r><w:t xml:space="preserve">If it </w:t></w:r><w:r><w:t xml:space="preserve">is not in this range, </w:t></w:r><w:r><w:t xml:space="preserve">a warning flag is set, </w:t></w:r><w:r><w:t>the torque is set to 0 and vehicle switches to Neutral State.</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>TorqueRequest_Nm</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is limited in the correct interval depending on the current </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>AutomaticTransmissionState</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> value: </w:t></w:r><w:r><w:t>if B between [-80, 80], if D between [0; 80], if R between [-40, 0], 0 if N or P.</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

if ($target -ne $null) {
    $null = $target.Range.InsertXML($xml)
} else {
    Write-Output "ERROR: anchor paragraph not found"
}
